$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44400
$ws.Range("M2").Value = 5
$ws.Range("N2").Value = 24000
$ws.Range("O2").Value = 24000
$ws.Range("P2").Value = 24000
$ws.Range("S2").Value = 1200

$ws.Range("D3").Value = 44419
$ws.Range("M3").Value = 40
$ws.Range("N3").Value = 25000
$ws.Range("O3").Value = 25000
$ws.Range("P3").Value = 25000
$ws.Range("S3").Value = 1250

$ws.Range("D4").Value = 44235
$ws.Range("M4").Value = 15
$ws.Range("N4").Value = 25000
$ws.Range("O4").Value = 25000
$ws.Range("P4").Value = 25000
$ws.Range("S4").Value = 1250

$ws.Range("D5").Value = 44412
$ws.Range("M5").Value = 20
$ws.Range("N5").Value = 25000
$ws.Range("O5").Value = 25000
$ws.Range("P5").Value = 25000
$ws.Range("S5").Value = 1250

$ws.Range("D6").Value = 44433
$ws.Range("M6").Value = 10
$ws.Range("N6").Value = 24000
$ws.Range("O6").Value = 24000
$ws.Range("P6").Value = 24000
$ws.Range("S6").Value = 1200

$ws.Range("D7").Value = 44334
$ws.Range("M7").Value = 20
$ws.Range("N7").Value = 25000
$ws.Range("O7").Value = 25000
$ws.Range("P7").Value = 25000
$ws.Range("S7").Value = 1250

$ws.Range("D8").Value = 44221
$ws.Range("M8").Value = 30
$ws.Range("N8").Value = 25000
$ws.Range("O8").Value = 25000
$ws.Range("P8").Value = 25000
$ws.Range("S8").Value = 1250

$ws.Range("D9").Value = 44421
$ws.Range("M9").Value = 20
$ws.Range("N9").Value = 24000
$ws.Range("O9").Value = 24000
$ws.Range("P9").Value = 24000
$ws.Range("S9").Value = 1200

$ws.Range("D10").Value = 44466
$ws.Range("M10").Value = 70
$ws.Range("N10").Value = 24000
$ws.Range("O10").Value = 24000
$ws.Range("P10").Value = 24000
$ws.Range("S10").Value = 1200

$ws.Range("D11").Value = 44462
$ws.Range("M11").Value = 10
$ws.Range("N11").Value = 24000
$ws.Range("O11").Value = 24000
$ws.Range("P11").Value = 24000
$ws.Range("S11").Value = 1200

$ws.Range("D12").Value = 44377
$ws.Range("M12").Value = 15
$ws.Range("N12").Value = 20000
$ws.Range("O12").Value = 20000
$ws.Range("P12").Value = 20000
$ws.Range("S12").Value = 1000

$ws.Range("D13").Value = 44435
$ws.Range("M13").Value = 100
$ws.Range("N13").Value = 24000
$ws.Range("O13").Value = 24000
$ws.Range("P13").Value = 24000
$ws.Range("S13").Value = 1200

$ws.Range("D14").Value = 44249
$ws.Range("M14").Value = 15
$ws.Range("N14").Value = 25000
$ws.Range("O14").Value = 25000
$ws.Range("P14").Value = 25000
$ws.Range("S14").Value = 1250

$ws.Range("D15").Value = 44356
$ws.Range("M15").Value = 15
$ws.Range("N15").Value = 24000
$ws.Range("O15").Value = 24000
$ws.Range("P15").Value = 24000
$ws.Range("S15").Value = 1200

$ws.Range("D16").Value = 44442
$ws.Range("M16").Value = 25
$ws.Range("N16").Value = 23000
$ws.Range("O16").Value = 23000
$ws.Range("P16").Value = 23000
$ws.Range("S16").Value = 1150

$ws.Range("D17").Value = 44363
$ws.Range("M17").Value = 30
$ws.Range("N17").Value = 24000
$ws.Range("O17").Value = 24000
$ws.Range("P17").Value = 24000
$ws.Range("S17").Value = 1200

$ws.Range("D18").Value = 44392
$ws.Range("M18").Value = 10
$ws.Range("N18").Value = 24000
$ws.Range("O18").Value = 24000
$ws.Range("P18").Value = 24000
$ws.Range("S18").Value = 1200

$ws.Range("D19").Value = 44426
$ws.Range("M19").Value = 15
$ws.Range("N19").Value = 24000
$ws.Range("O19").Value = 24000
$ws.Range("P19").Value = 24000
$ws.Range("S19").Value = 1200

$ws.Range("D20").Value = 44382
$ws.Range("M20").Value = 15
$ws.Range("N20").Value = 20000
$ws.Range("O20").Value = 20000
$ws.Range("P20").Value = 20000
$ws.Range("S20").Value = 1000

$ws.Range("D21").Value = 44431
$ws.Range("M21").Value = 40
$ws.Range("N21").Value = 24000
$ws.Range("O21").Value = 24000
$ws.Range("P21").Value = 24000
$ws.Range("S21").Value = 1200

$ws.Range("D22").Value = 44424
$ws.Range("M22").Value = 25
$ws.Range("N22").Value = 24000
$ws.Range("O22").Value = 24000
$ws.Range("P22").Value = 24000
$ws.Range("S22").Value = 1200

$ws.Range("D23").Value = 44488
$ws.Range("M23").Value = 40
$ws.Range("N23").Value = 20000
$ws.Range("O23").Value = 20000
$ws.Range("P23").Value = 20000
$ws.Range("S23").Value = 1000

$ws.Range("D24").Value = 44214
$ws.Range("M24").Value = 15
$ws.Range("N24").Value = 25000
$ws.Range("O24").Value = 25000
$ws.Range("P24").Value = 25000
$ws.Range("S24").Value = 1250

$ws.Range("D25").Value = 44349
$ws.Range("M25").Value = 30
$ws.Range("N25").Value = 24000
$ws.Range("O25").Value = 24000
$ws.Range("P25").Value = 24000
$ws.Range("S25").Value = 1200

$ws.Range("D26").Value = 44222
$ws.Range("M26").Value = 15
$ws.Range("N26").Value = 25000
$ws.Range("O26").Value = 25000
$ws.Range("P26").Value = 25000
$ws.Range("S26").Value = 1250

$ws.Range("D27").Value = 44390
$ws.Range("M27").Value = 10
$ws.Range("N27").Value = 24000
$ws.Range("O27").Value = 24000
$ws.Range("P27").Value = 24000
$ws.Range("S27").Value = 1200

$ws.Range("D28").Value = 44469
$ws.Range("M28").Value = 40
$ws.Range("N28").Value = 24000
$ws.Range("O28").Value = 24000
$ws.Range("P28").Value = 24000
$ws.Range("S28").Value = 1200

$ws.Range("D29").Value = 44475
$ws.Range("M29").Value = 20
$ws.Range("N29").Value = 24000
$ws.Range("O29").Value = 24000
$ws.Range("P29").Value = 24000
$ws.Range("S29").Value = 1200

$ws.Range("D30").Value = 44468
$ws.Range("M30").Value = 20
$ws.Range("N30").Value = 24000
$ws.Range("O30").Value = 24000
$ws.Range("P30").Value = 24000
$ws.Range("S30").Value = 1200

$ws.Range("D31").Value = 44532
$ws.Range("M31").Value = 20
$ws.Range("N31").Value = 28000
$ws.Range("O31").Value = 28000
$ws.Range("P31").Value = 28000
$ws.Range("S31").Value = 1400

$ws.Range("D32").Value = 44489
$ws.Range("M32").Value = 40
$ws.Range("N32").Value = 24000
$ws.Range("O32").Value = 24000
$ws.Range("P32").Value = 24000
$ws.Range("S32").Value = 1200

$ws.Range("D33").Value = 44425
$ws.Range("M33").Value = 15
$ws.Range("N33").Value = 24000
$ws.Range("O33").Value = 24000
$ws.Range("P33").Value = 24000
$ws.Range("S33").Value = 1200

$ws.Range("D34").Value = 44418
$ws.Range("M34").Value = 20
$ws.Range("N34").Value = 24000
$ws.Range("O34").Value = 24000
$ws.Range("P34").Value = 24000
$ws.Range("S34").Value = 1200

$ws.Range("D35").Value = 44434
$ws.Range("M35").Value = 20
$ws.Range("N35").Value = 24000
$ws.Range("O35").Value = 24000
$ws.Range("P35").Value = 24000
$ws.Range("S35").Value = 1200

$ws.Range("D36").Value = 44467
$ws.Range("M36").Value = 20
$ws.Range("N36").Value = 24000
$ws.Range("O36").Value = 24000
$ws.Range("P36").Value = 24000
$ws.Range("S36").Value = 1200

$ws.Range("D37").Value = 44231
$ws.Range("M37").Value = 15
$ws.Range("N37").Value = 25000
$ws.Range("O37").Value = 25000
$ws.Range("P37").Value = 25000
$ws.Range("S37").Value = 1250

$ws.Range("D38").Value = 44251
$ws.Range("M38").Value = 15
$ws.Range("N38").Value = 25000
$ws.Range("O38").Value = 25000
$ws.Range("P38").Value = 25000
$ws.Range("S38").Value = 1250

$ws.Range("D39").Value = 44389
$ws.Range("M39").Value = 20
$ws.Range("N39").Value = 24000
$ws.Range("O39").Value = 24000
$ws.Range("P39").Value = 24000
$ws.Range("S39").Value = 1200

$ws.Range("D40").Value = 44391
$ws.Range("M40").Value = 10
$ws.Range("N40").Value = 24000
$ws.Range("O40").Value = 24000
$ws.Range("P40").Value = 24000
$ws.Range("S40").Value = 1200

$ws.Range("D41").Value = 44396
$ws.Range("M41").Value = 12
$ws.Range("N41").Value = 24000
$ws.Range("O41").Value = 24000
$ws.Range("P41").Value = 24000
$ws.Range("S41").Value = 1200

$ws.Range("D42").Value = 44232
$ws.Range("M42").Value = 15
$ws.Range("N42").Value = 25000
$ws.Range("O42").Value = 25000
$ws.Range("P42").Value = 25000
$ws.Range("S42").Value = 1250

$ws.Range("D43").Value = 44420
$ws.Range("M43").Value = 35
$ws.Range("N43").Value = 25000
$ws.Range("O43").Value = 25000
$ws.Range("P43").Value = 25000
$ws.Range("S43").Value = 1250

$ws.Range("D44").Value = 44414
$ws.Range("M44").Value = 15
$ws.Range("N44").Value = 25000
$ws.Range("O44").Value = 25000
$ws.Range("P44").Value = 25000
$ws.Range("S44").Value = 1250

$ws.Range("D45").Value = 44474
$ws.Range("M45").Value = 20
$ws.Range("N45").Value = 24000
$ws.Range("O45").Value = 24000
$ws.Range("P45").Value = 24000
$ws.Range("S45").Value = 1200

$ws.Range("D46").Value = 44238
$ws.Range("M46").Value = 30
$ws.Range("N46").Value = 25000
$ws.Range("O46").Value = 25000
$ws.Range("P46").Value = 25000
$ws.Range("S46").Value = 1250

$ws.Range("D47").Value = 44175
$ws.Range("M47").Value = 25
$ws.Range("N47").Value = 23000
$ws.Range("O47").Value = 23000
$ws.Range("P47").Value = 23000
$ws.Range("S47").Value = 1150

$ws.Range("D48").Value = 44432
$ws.Range("M48").Value = 30
$ws.Range("N48").Value = 24000
$ws.Range("O48").Value = 24000
$ws.Range("P48").Value = 24000
$ws.Range("S48").Value = 1200

$ws.Range("D49").Value = 44428
$ws.Range("M49").Value = 15
$ws.Range("N49").Value = 24000
$ws.Range("O49").Value = 24000
$ws.Range("P49").Value = 24000
$ws.Range("S49").Value = 1200

$ws.Range("D50").Value = 44452
$ws.Range("M50").Value = 25
$ws.Range("N50").Value = 25000
$ws.Range("O50").Value = 25000
$ws.Range("P50").Value = 25000
$ws.Range("S50").Value = 1250

$ws.Range("D51").Value = 44398
$ws.Range("M51").Value = 15
$ws.Range("N51").Value = 25000
$ws.Range("O51").Value = 25000
$ws.Range("P51").Value = 25000
$ws.Range("S51").Value = 1250

$ws.Range("D52").Value = 44461
$ws.Range("M52").Value = 30
$ws.Range("N52").Value = 24000
$ws.Range("O52").Value = 24000
$ws.Range("P52").Value = 24000
$ws.Range("S52").Value = 1200

$ws.Range("D53").Value = 44454
$ws.Range("M53").Value = 25
$ws.Range("N53").Value = 25000
$ws.Range("O53").Value = 25000
$ws.Range("P53").Value = 25000
$ws.Range("S53").Value = 1250

$ws.Range("D54").Value = 44194
$ws.Range("M54").Value = 20
$ws.Range("N54").Value = 20000
$ws.Range("O54").Value = 20000
$ws.Range("P54").Value = 20000
$ws.Range("S54").Value = 1000
